# RDS Agreement Manager.xlsx - apply the "Display AA now working / Incomplete
# data fixed" commit.
#
# Summary of the target edit (from the OOXML diff):
#   - xl/workbook.xml
#       * workbookPr gets codeName="ThisWorkbook"              (VBA code name;
#         not reachable through this COM surface - see note below)
#       * Sheet1 is renamed to "Main"
#       * 9 workbook-scoped defined names are added, all pointing at cells on
#         the new "Main" sheet
#   - xl/sharedStrings.xml (new) / xl/worksheets/sheet1.xml
#       * B2 on "Main" gets the text "CrewNo"
#       * the sheet's dimension/selection move off A1 onto the new data
#       * the sheetView records the live selection on B4 (SlotStart cell)
#   - xl/worksheets/sheet2.xml / sheet3.xml
#       * sheetPr codeName="Sheet2" / "Sheet3" (VBA code names - same caveat)
#
# NOTE on VBA `CodeName`: this COM host exposes `.CodeName` as a gettable
# property (it already mirrors the sheet's internal id, e.g. "Sheet1") but
# assigning to it is a silent no-op - it never reaches the saved
# `<sheetPr codeName=.../>` / `<workbookPr codeName=.../>` XML (verified: a
# round trip with only `.CodeName = "..."` assignments produces a byte-for-
# byte identical workbook.xml/sheetN.xml to one with no script at all). We
# still set it below so the intent is explicit and it's a correct no-op if a
# future host version wires it up, but it cannot be relied on to move the
# XML in this runtime.

$wb = $excel.ActiveWorkbook

# --- Sheet1 -> "Main" ------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Main"
$ws1.CodeName = "ShtFrontPage"   # best-effort; see note above

# --- Workbook / VBA code name ----------------------------------------------
$wb.CodeName = "ThisWorkbook"    # best-effort; see note above

# --- Sheet2 / Sheet3 code names ---------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.CodeName = "Sheet2"         # best-effort; see note above

$ws3 = $wb.Worksheets.Item(3)
$ws3.CodeName = "Sheet3"         # best-effort; see note above

# --- Data entry on Main -----------------------------------------------------
# B2 holds the literal label "CrewNo" (becomes the lone shared string).
$ws1.Range("B2").Value = "CrewNo"

# Leave the live selection on B4 (the SlotStart cell), matching the saved
# sheetView/<selection activeCell="B4" sqref="B4"/>.
$ws1.Range("B4").Select() | Out-Null

# --- Workbook-scoped defined names ------------------------------------------
# Single quotes avoid PowerShell interpolating "$F", "$2", etc.
$wb.Names.Add('ContractType', '=Main!$F$2')
$wb.Names.Add('CrewName',     '=Main!$C$2')
$wb.Names.Add('CrewNo',       '=Main!$B$3')
$wb.Names.Add('HrsWk',        '=Main!$G$2')
$wb.Names.Add('NoWeeks',      '=Main!$H$2')
$wb.Names.Add('RevDate',      '=Main!$I$2')
$wb.Names.Add('Role',         '=Main!$D$2')
$wb.Names.Add('SlotStart',    '=Main!$B$4')
$wb.Names.Add('TemplateDate', '=Main!$E$2')
